$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value (Excel auto-converts numeric-looking
# strings like "506.60" or "1.00" into real numbers on assignment, just like
# typing them into a General-formatted cell). Temporarily mark the cell as
# Text, assign the literal string, then restore the default "Normal" style
# so no stray number-format is left behind.
function Set-TextValue {
    param($Address, $Text)
    $range = $ws.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

# --- Plain value updates (unambiguous text: URLs, names, percent strings, ---
# --- and prices using "." as a thousands separator so they cannot parse ---
# --- as a single number) ---
$ws.Range("D2").Value = "57.013.65"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.402.45"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "2.412.49"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").Value = "2.829.49"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "56.948.70"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +2.62%  "
$ws.Range("D18").Value = "2.433.28"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "0.0₃0728"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("E40").Value = "  +5.83%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E41").Value = "  +3.33%  "
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("E51").Value = "  +9.09%  "

# --- Price values that look like plain decimal numbers and would otherwise ---
# --- be silently converted from text to a number by Excel ---
Set-TextValue "D5" "506.60"
Set-TextValue "D6" "132.86"
Set-TextValue "D8" "0.554"
Set-TextValue "D10" "0.0978"
Set-TextValue "D12" "0.323"
Set-TextValue "D13" "4.59"
Set-TextValue "D16" "21.83"
Set-TextValue "D17" "0.0000135"
Set-TextValue "D19" "10.27"
Set-TextValue "D20" "4.06"
Set-TextValue "D21" "310.39"
Set-TextValue "D22" "6.28"
Set-TextValue "D23" "1.00"
Set-TextValue "D24" "5.63"
Set-TextValue "D25" "67.01"
Set-TextValue "D27" "0.378"
Set-TextValue "D28" "0.153"
Set-TextValue "D29" "7.44"
Set-TextValue "D30" "175.79"
Set-TextValue "D32" "1.68"
Set-TextValue "D34" "5.91"
Set-TextValue "D37" "17.99"
Set-TextValue "D38" "1.20"
Set-TextValue "D39" "3.83"
Set-TextValue "D40" "0.830"
Set-TextValue "D41" "36.85"
Set-TextValue "D42" "1.45"
Set-TextValue "D43" "133.09"
Set-TextValue "D44" "3.38"
Set-TextValue "D45" "4.88"
Set-TextValue "D46" "0.570"
Set-TextValue "D47" "251.64"
Set-TextValue "D48" "0.0913"
Set-TextValue "D49" "0.0490"
Set-TextValue "D51" "17.21"
